# feat: change language to russian
# (the underlying OOXML diff only appends three new delivery-note rows;
#  this script reproduces exactly that data change)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value as TEXT (matches the source
# workbook's convention of storing Price/TotalAmount as inline strings)
# without leaving a stray custom number-format style behind.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$rows = @(
    @{ Row = 8;  A = 18; B = "Product 7";  C = 4; D = "60000.00";   E = 45681.37577218791; F = 9;  G = "240000.00" },
    @{ Row = 9;  A = 19; B = "fesfesfes";  C = 5; D = "343434.00";  E = 45681.40553633918; F = 11; G = "1717170.00" },
    @{ Row = 10; A = 20; B = "Product 8";  C = 1; D = "12430.00";   E = 45681.42420941254; F = 10; G = "12430.00" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C

    Set-TextValue $row 4 $r.D

    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 5).NumberFormat = $ws.Cells.Item(7, 5).NumberFormat

    $ws.Cells.Item($row, 6).Value = $r.F

    Set-TextValue $row 7 $r.G
}

Write-Output "Appended rows 8-10 to Sheet1"
